$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.475.74"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "'3.550.85"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'203.64"
$ws.Range("E5").Value = "  +5.79%  "
$ws.Range("D6").Value = "'555.82"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").Value = "'3.528.41"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D10").Value = "'0.663"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'61.35"
$ws.Range("E11").Value = "  +11.36%  "
$ws.Range("D12").Value = "'0.145"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "'9.99"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'4.109.89"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "'3.533.85"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "'18.77"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "'67.186.88"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'393.30"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "'12.23"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'4.05"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "'83.23"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").Value = "'3.80"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").Value = "'12.15"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "'8.98"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D32").Value = "'7.33"
$ws.Range("E32").Value = "  -6.53%  "
$ws.Range("D33").Value = "'11.87"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "'63.49"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'0.112"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "'40.67"
$ws.Range("E36").Value = "  -3.82%  "
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'3.10"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D45").Value = "'2.57"
$ws.Range("E45").Value = "  -6.85%  "
$ws.Range("E46").Value = "  +7.45%  "
$ws.Range("D47").Value = "'0.0403"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "'0.128"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").Value = "'2.97"
$ws.Range("E49").Value = "  -3.57%  "
$ws.Range("B30").Value = "'Bittensor"
$ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'711.01"
$ws.Range("E30").Value = "  +9.27%  "
$ws.Range("B31").Value = "'EthereumClassic"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'30.96"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("B40").Value = "'Maker"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'3.122.20"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("B41").Value = "'PEPE"
$ws.Range("C41").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0716"
$ws.Range("E41").Value = "  -6.48%  "
$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("B44").Value = "'dogwifhat"
$ws.Range("C44").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.94"
$ws.Range("E44").Value = "  +21.74%  "
$ws.Range("B50").Value = "'Monero"
$ws.Range("C50").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'136.85"
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("B51").Value = "'THORChain"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.37"
$ws.Range("E51").Value = "  -1.78%  "
